$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Insert a new row above row 12 ("Programa resumido:" block), shifting everything
# from row 12 downward by one row. This makes room for the new
# "Docentes responsáveis:" label row. The inserted row already carries over the
# A-column bold style, and all subsequent row heights line up with the target
# layout without any further adjustment.
$ws.Rows(12).Insert()

# New label-only row 12: "Docentes responsáveis:"
$ws.Range("A12").Value = 'Docentes responsáveis:'

# Row 13 ("Programa resumido:") - label unchanged, data cells now hold the
# updated activation date text. Excel auto-detects the date-shaped string and
# would store it as a date serial with a new number-formatted style, so we
# force it back to plain text matching the existing column B/C styles.
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '01/01/2023'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = '01/01/2023'
$ws.Range("B3").Copy()
$ws.Range("B13").PasteSpecial($xlPasteFormats)
$ws.Range("C3").Copy()
$ws.Range("C13").PasteSpecial($xlPasteFormats)

# Row 14 ("Short syllabus:") now carries the short syllabus text directly in
# brand-new B14/C14 cells (they did not exist before). Freshly created cells
# pick up the bold A-column style by default, so re-apply the correct
# column B/C (plain / red) styles afterward.
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Review of Newtonian dynamics (“vector mechanics”). General structure of classical mechanics in the Lagrangian and Hamiltonian formulations. Applications to problems of central forces and dynamics of rigid bodies. Problems in non-inertial frames.'
$ws.Range("C14").Value = 'Review of Newtonian dynamics (“vector mechanics”). General structure of classical mechanics in the Lagrangian and Hamiltonian formulations. Applications to problems of central forces and dynamics of rigid bodies. Problems in non-inertial frames.'
$ws.Range("B3").Copy()
$ws.Range("B14").PasteSpecial($xlPasteFormats)
$ws.Range("C3").Copy()
$ws.Range("C14").PasteSpecial($xlPasteFormats)

# Row 15 ("Programa:") now shows the responsible professor.
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Range("C15").Value = '5840730 - Antonio Jefferson da Silva Machado'

# Row 16 ("Syllabus:") now carries the full syllabus text directly in
# brand-new B16/C16 cells; fix their style the same way as row 14.
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = "1. Review of point kinematics: position, velocity and acceleration vectors. 2. Forces, resultant force, conservation of linear momentum and Newton's Laws of Dynamics; 3. Work and energy; kinetic energy and potential energy. Kinetic energy theorem and conservation of total energy. 4. Hamilton's variational principle and the Lagrangian Formulation of Mechanics: generalized coordinates, Lagrangian and the Euler-Lagrange equations. Ignorable coordinates and conservation principles. 5. Applications: central forces; torque and conservation of angular momentum; two-body problems with mutual attraction or repulsion; Gravitation and Kepler's laws. 6. Linearly accelerated and rotating non-inertial frames of reference. Inertia forces: centrifugal force, Coriolis force. Effects of inertia forces on planet Earth; Foucault pendulum. 7. Statics and Dynamics of Rigid Bodies; moments of inertia; torques; Plane motion of rigid bodies; parallel axis theorem. 8. Hamilton's Formulation for Classical Mechanics: the Hamiltonian and Hamilton's equations"
$ws.Range("C16").Value = "1. Review of point kinematics: position, velocity and acceleration vectors. 2. Forces, resultant force, conservation of linear momentum and Newton's Laws of Dynamics; 3. Work and energy; kinetic energy and potential energy. Kinetic energy theorem and conservation of total energy. 4. Hamilton's variational principle and the Lagrangian Formulation of Mechanics: generalized coordinates, Lagrangian and the Euler-Lagrange equations. Ignorable coordinates and conservation principles. 5. Applications: central forces; torque and conservation of angular momentum; two-body problems with mutual attraction or repulsion; Gravitation and Kepler's laws. 6. Linearly accelerated and rotating non-inertial frames of reference. Inertia forces: centrifugal force, Coriolis force. Effects of inertia forces on planet Earth; Foucault pendulum. 7. Statics and Dynamics of Rigid Bodies; moments of inertia; torques; Plane motion of rigid bodies; parallel axis theorem. 8. Hamilton's Formulation for Classical Mechanics: the Hamiltonian and Hamilton's equations"
$ws.Range("B3").Copy()
$ws.Range("B16").PasteSpecial($xlPasteFormats)
$ws.Range("C3").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)

# Row 18 ("Método:") - new responsible professor for the grading method.
$ws.Range("B18").Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Range("C18").Value = '1176388 - Luiz Tadeu Fernandes Eleno'

# Row 20 ("Norma de recuperação:") - updated make-up exam rule text.
$ws.Range("B20").Value = 'A nota final (NF) será a média ponderada de três provas, P1 (peso 1), P2 (peso 1) e P3 (peso 2)'
$ws.Range("C20").Value = 'A nota final (NF) será a média ponderada de três provas, P1 (peso 1), P2 (peso 1) e P3 (peso 2)'
